# "Ran Test on Web App" -- mark every data row on every sheet with a
# "Test Result" = "Success" value in the first empty column, and leave
# the selection/active-sheet state the way it was left after testing
# each sheet in turn (Customers -> Order Details -> Orders -> Products).

$wb = $excel.ActiveWorkbook

# --- Customers --------------------------------------------------------
$wsCustomers = $wb.Worksheets.Item("Customers")
$wsCustomers.Range("F2:F14").Value = "Success"
$wsCustomers.Range("K27").Select() | Out-Null

# --- Order Details ------------------------------------------------------
$wsOrderDetails = $wb.Worksheets.Item("Order Details")
$wsOrderDetails.Range("F2:F41").Value = "Success"
$wsOrderDetails.Range("G34").Select() | Out-Null

# --- Orders -------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Orders")
$wsOrders.Range("E2:E15").Value = "Success"

# --- Products -------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Range("E2:E40").Value = "Success"
$wsProducts.Range("H19").Select() | Out-Null
